$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows for Group 205 / Script "05" ------------------------------
$data = @(
    @(205, "05", "drumming",          8,  9, 3),
    @(205, "05", "clicking pen",      8,  7, 6),
    @(205, "05", "head on table",     2,  8, 6),
    @(205, "05", "locking at phone", 10,  9, 7),
    @(205, "05", "snipping",          6,  9, 7),
    @(205, "05", "drawing",           1,  9, 9),
    @(205, "05", "heckling",          9,  3, 5),
    @(205, "05", "whispering",        7,  6, 7),
    @(205, "05", "chatting",          5,  9, 4)
)

$startRow = 38
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# --- Header row: rename/reorder the factor columns -------------------------
# Old: D1=Disturbing_Factor, E1=Confident_Factor, F1=Dispersion_School
# New: D1=Disruption_Factor, E1=Confident_Factor, F1=Dispersion_School
$ws.Range("D1").Value = "Disruption_Factor"
$ws.Range("E1").Value = "Confident_Factor"
$ws.Range("F1").Value = "Dispersion_School"

# --- Selection change --------------------------------------------------------
$ws.Range("D1").Select()
